$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''43.885.74'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +4.55%  '
$ws.Range('D3').Value = '''2.281.99'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +2.41%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''231.73'
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Value = '''0.629'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.67%  '
$ws.Range('D7').Value = '''61.70'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.91%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '''0.422'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +5.02%  '
$ws.Range('D10').Value = '''0.0947'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +6.42%  '
$ws.Range('D11').Value = '''57.91'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.80%  '
$ws.Range('E12').Value = '  +0.66%  '
$ws.Range('D13').Value = '''2.620.88'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +2.43%  '
$ws.Range('D14').Value = '''15.84'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.23%  '
$ws.Range('D15').Value = '''23.75'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +9.21%  '
$ws.Range('D16').Value = '''5.83'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +4.73%  '
$ws.Range('D17').Value = '''0.814'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +2.01%  '
$ws.Range('D18').Value = '''2.280.95'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +2.28%  '
$ws.Range('D19').Value = '''43.760.48'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +4.60%  '
$ws.Range('D20').Value = '''0.0₃0939'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +5.47%  '
$ws.Range('D21').Value = '''73.24'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.99%  '
$ws.Range('D22').Value = '''6.26'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +3.90%  '
$ws.Range('D23').Value = '''251.14'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.46%  '
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').Value = '''2.56'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +7.63%  '
$ws.Range('E26').Value = '  +2.42%  '
$ws.Range('D27').Value = '''9.89'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.16%  '
$ws.Range('D28').Value = '''171.16'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +2.25%  '
$ws.Range('E29').Value = '  -0.64%  '
$ws.Range('D30').Value = '''20.63'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +3.52%  '
$ws.Range('E31').Value = '  +4.72%  '
$ws.Range('D32').Value = '''2.66'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.07%  '
$ws.Range('E33').Value = '  +0.71%  '
$ws.Range('D34').Value = '''4.81'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +4.21%  '
$ws.Range('D35').Value = '''5.06'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +2.44%  '
$ws.Range('E36').Value = '  +5.53%  '
$ws.Range('D37').Value = '''6.52'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.80%  '
$ws.Range('E38').Value = '  +2.86%  '
$ws.Range('E39').Value = '  -1.43%  '
$ws.Range('E40').Value = '  +4.79%  '
$ws.Range('D41').Value = '''0.998'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.36%  '
$ws.Range('D42').Value = '''8.78'
$ws.Range('D42').Style = "Normal"
$ws.Range('D43').Value = '''0.000226'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -12.11%  '
$ws.Range('D44').Value = '''4.54'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -5.26%  '
$ws.Range('D45').Value = '''0.0981'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('E46').Value = '  +0.81%  '
$ws.Range('D47').Value = '''98.35'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.53%  '
$ws.Range('D48').Value = '''1.472.17'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.10%  '
$ws.Range('D49').Value = '''16.72'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.48%  '
$ws.Range('E50').Value = '  +1.24%  '
$ws.Range('E51').Value = '  -1.13%  '
